# Applies numeric value corrections to the Chocobo Profits leve-crafting
# workbook, one worksheet/row at a time, matching the source diff.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 25999.5
$ws.Range("J7").Value = 25999.5
$ws.Range("L7").Value = 25999.5
$ws.Range("N7").Value = -26223.5

# Row 10
$ws.Range("H10").Value = 29499.75
$ws.Range("J10").Value = 29499.75
$ws.Range("L10").Value = 29499.75
$ws.Range("N10").Value = -30085.75

# Row 14
$ws.Range("H14").Value = 25999.5
$ws.Range("J14").Value = 25999.5
$ws.Range("L14").Value = 25999.5
$ws.Range("N14").Value = -26381.5

# Row 33
$ws.Range("H33").Value = 350.93332
$ws.Range("I33").Value = 250.92308
$ws.Range("K33").Value = 250.92308
$ws.Range("M33").Value = -21.92308

# Row 74
$ws.Range("H74").Value = 4766563
$ws.Range("I74").Value = 8337283.5
$ws.Range("J74").Value = 5602.222
$ws.Range("K74").Value = 8337283.5
$ws.Range("L74").Value = 5602.222
$ws.Range("M74").Value = -8336347.5
$ws.Range("N74").Value = -7474.222

# Row 77
$ws.Range("H77").Value = 4766563
$ws.Range("I77").Value = 8337283.5
$ws.Range("J77").Value = 5602.222
$ws.Range("K77").Value = 41686417.5
$ws.Range("L77").Value = 28011.11
$ws.Range("M77").Value = -41681737.5
$ws.Range("N77").Value = -37371.11

# Row 99
$ws.Range("H99").Value = 972
$ws.Range("I99").Value = 664
$ws.Range("J99").Value = 1280
$ws.Range("K99").Value = 1992
$ws.Range("L99").Value = 3840
$ws.Range("M99").Value = -494
$ws.Range("N99").Value = -6836

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5331.5957
$ws.Range("I32").Value = 4660.436
$ws.Range("J32").Value = 8603.5
$ws.Range("K32").Value = 4660.436
$ws.Range("L32").Value = 8603.5
$ws.Range("M32").Value = -4373.436
$ws.Range("N32").Value = -9177.5

# Row 61
$ws.Range("H61").Value = 2138.4
$ws.Range("I61").Value = 927
$ws.Range("K61").Value = 927
$ws.Range("M61").Value = -715

# Row 74
$ws.Range("H74").Value = 7171.4375
$ws.Range("I74").Value = 8222.637000000001
$ws.Range("J74").Value = 4858.8
$ws.Range("K74").Value = 8222.637000000001
$ws.Range("L74").Value = 4858.8
$ws.Range("M74").Value = -7348.637000000001
$ws.Range("N74").Value = -6606.8

# Row 77
$ws.Range("H77").Value = 7171.4375
$ws.Range("I77").Value = 8222.637000000001
$ws.Range("J77").Value = 4858.8
$ws.Range("K77").Value = 41113.185
$ws.Range("L77").Value = 24294
$ws.Range("M77").Value = -36745.185
$ws.Range("N77").Value = -33030

# Row 97
$ws.Range("H97").Value = 1721.1875
$ws.Range("I97").Value = 1180.5714
$ws.Range("K97").Value = 1180.5714
$ws.Range("M97").Value = -684.5714

# Row 136
$ws.Range("H136").Value = 2138.4
$ws.Range("I136").Value = 927
$ws.Range("K136").Value = 2781
$ws.Range("M136").Value = -231

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1122.8889
$ws.Range("I94").Value = 888.25
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 888.25
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -437.25
$ws.Range("N94").Value = -3902

# Row 134
$ws.Range("H134").Value = 1964.8182
$ws.Range("I134").Value = 1287.8214
$ws.Range("J134").Value = 5756
$ws.Range("K134").Value = 3863.4642
$ws.Range("L134").Value = 17268
$ws.Range("M134").Value = -1328.4642
$ws.Range("N134").Value = -22338

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3225.4167
$ws.Range("I31").Value = 1113.1333
$ws.Range("J31").Value = 6745.8887
$ws.Range("K31").Value = 1113.1333
$ws.Range("L31").Value = 6745.8887
$ws.Range("M31").Value = -818.1333
$ws.Range("N31").Value = -7335.8887

# Row 34
$ws.Range("H34").Value = 3225.4167
$ws.Range("I34").Value = 1113.1333
$ws.Range("J34").Value = 6745.8887
$ws.Range("K34").Value = 1113.1333
$ws.Range("L34").Value = 6745.8887
$ws.Range("M34").Value = -911.1333
$ws.Range("N34").Value = -7149.8887

# Row 58
$ws.Range("H58").Value = 2966.1167
$ws.Range("I58").Value = 1724.8431
$ws.Range("K58").Value = 1724.8431
$ws.Range("M58").Value = -1521.8431

# Row 136
$ws.Range("H136").Value = 2966.1167
$ws.Range("I136").Value = 1724.8431
$ws.Range("K136").Value = 5174.5293
$ws.Range("M136").Value = -2624.5293

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 728.0599999999999
$ws.Range("I131").Value = 321.13333
$ws.Range("J131").Value = 799.8706
$ws.Range("K131").Value = 963.39999
$ws.Range("L131").Value = 2399.6118
$ws.Range("M131").Value = 4076.60001
$ws.Range("N131").Value = -12479.6118

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 3505530.8
$ws.Range("I3").Value = 8752000
$ws.Range("J3").Value = 7884.5
$ws.Range("K3").Value = 8752000
$ws.Range("L3").Value = 7884.5
$ws.Range("M3").Value = -8751884
$ws.Range("N3").Value = -8116.5

# Row 10
$ws.Range("H10").Value = 9001680
$ws.Range("I10").Value = 11251250
$ws.Range("J10").Value = 3400
$ws.Range("K10").Value = 11251250
$ws.Range("L10").Value = 3400
$ws.Range("M10").Value = -11251081
$ws.Range("N10").Value = -3738

# Row 64
$ws.Range("H64").Value = 32676.4
$ws.Range("J64").Value = 32676.4
$ws.Range("L64").Value = 32676.4
$ws.Range("N64").Value = -33172.4

# Row 67
$ws.Range("H67").Value = 32676.4
$ws.Range("J67").Value = 32676.4
$ws.Range("L67").Value = 32676.4
$ws.Range("N67").Value = -34392.4

# Row 97
$ws.Range("H97").Value = 934.95
$ws.Range("I97").Value = 938.0625
$ws.Range("J97").Value = 922.5
$ws.Range("K97").Value = 938.0625
$ws.Range("L97").Value = 922.5
$ws.Range("M97").Value = -442.0625
$ws.Range("N97").Value = -1914.5

# Row 132
$ws.Range("H132").Value = 6928.4287
$ws.Range("J132").Value = 7083.1665
$ws.Range("L132").Value = 21249.4995
$ws.Range("N132").Value = -26309.4995

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 96
$ws.Range("H96").Value = 30197
$ws.Range("J96").Value = 30197
$ws.Range("L96").Value = 30197
$ws.Range("N96").Value = -35689

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 703
$ws.Range("I6").Value = 703
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 703
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -588
$ws.Range("N6").ClearContents()

# Row 9
$ws.Range("H9").Value = 5980
$ws.Range("J9").Value = 5980
$ws.Range("L9").Value = 5980
$ws.Range("N9").Value = -6260

# Row 96
$ws.Range("H96").Value = 142931920
$ws.Range("I96").Value = 200100700
$ws.Range("J96").Value = 9950
$ws.Range("K96").Value = 200100700
$ws.Range("L96").Value = 9950
$ws.Range("M96").Value = -200099327
$ws.Range("N96").Value = -12696

# Row 133
$ws.Range("H133").Value = 41555.832
$ws.Range("J133").Value = 41555.832
$ws.Range("L133").Value = 41555.832
$ws.Range("N133").Value = -51675.832

# Row 136
$ws.Range("H136").Value = 3841.7932
$ws.Range("I136").Value = 1881.2858
$ws.Range("K136").Value = 5643.857400000001
$ws.Range("M136").Value = -3093.857400000001
